$wb = $excel.ActiveWorkbook

# zh-cn sheet: update handoff/handback datetime for the 8aae... row (row 2)
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("D2").Value = "2016-02-06 03:59:21"
$ws.Range("G2").Value = "2016-02-06 04:00:06"

# de-de sheet: update handoff/handback datetime for the 8aae... row (row 2)
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("D2").Value = "2016-02-06 03:59:32"
$ws.Range("G2").Value = "2016-02-06 04:00:25"
